$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 98999.8
$ws.Range("J3").Value = 98999.8
$ws.Range("L3").Value = 98999.8
$ws.Range("N3").Value = -99227.8
$ws.Range("H12").Value = 368.6
$ws.Range("I12").Value = 247.66667
$ws.Range("K12").Value = 247.66667
$ws.Range("M12").Value = -77.66667000000001
$ws.Range("H99").Value = 3037.7
$ws.Range("J99").Value = 4955.8335
$ws.Range("L99").Value = 14867.5005
$ws.Range("N99").Value = -17863.5005
$ws.Range("H102").Value = 98999.8
$ws.Range("J102").Value = 98999.8
$ws.Range("L102").Value = 98999.8
$ws.Range("N102").Value = -105489.8
$ws.Range("H112").Value = 2008.625
$ws.Range("J112").Value = 1964.7142
$ws.Range("L112").Value = 5894.142599999999
$ws.Range("N112").Value = -8110.142599999999
$ws.Range("H118").Value = 963.619
$ws.Range("I118").Value = 501
$ws.Range("K118").Value = 1503
$ws.Range("M118").Value = 154
$ws.Range("H127").Value = 2698.6
$ws.Range("I127").Value = 2548.25
$ws.Range("K127").Value = 7644.75
$ws.Range("M127").Value = -2684.75
$ws.Range("H129").Value = 889.5833
$ws.Range("I129").Value = 889.5833
$ws.Range("K129").Value = 2668.7499
$ws.Range("M129").Value = 2331.2501
$ws.Range("H131").Value = 2951.2727
$ws.Range("I131").Value = 1246.4
$ws.Range("J131").Value = 20000
$ws.Range("K131").Value = 3739.2
$ws.Range("L131").Value = 60000
$ws.Range("M131").Value = 1300.8
$ws.Range("N131").Value = -70080
$ws.Range("H138").Value = 2508.9016
$ws.Range("I138").Value = 2275.96
$ws.Range("J138").Value = 2670.6667
$ws.Range("K138").Value = 6827.88
$ws.Range("L138").Value = 8012.000100000001
$ws.Range("M138").Value = -1687.88
$ws.Range("N138").Value = -18292.0001
$ws.Range("H141").Value = 4956.1177
$ws.Range("I141").Value = 5430.3335
$ws.Range("K141").Value = 16291.0005
$ws.Range("M141").Value = -11111.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2869.4167
$ws.Range("I45").Value = 1037.0834
$ws.Range("J45").Value = 4701.75
$ws.Range("K45").Value = 1037.0834
$ws.Range("L45").Value = 4701.75
$ws.Range("M45").Value = -660.0834
$ws.Range("N45").Value = -5455.75
$ws.Range("H132").Value = 17053.54
$ws.Range("I132").Value = 21832.652
$ws.Range("J132").Value = 2417.5
$ws.Range("K132").Value = 65497.95599999999
$ws.Range("L132").Value = 7252.5
$ws.Range("M132").Value = -62967.95599999999
$ws.Range("N132").Value = -12312.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1194.1333
$ws.Range("I94").Value = 1084.1852
$ws.Range("K94").Value = 1084.1852
$ws.Range("M94").Value = -633.1851999999999
$ws.Range("H134").Value = 2524.6428
$ws.Range("I134").Value = 1948.5385
$ws.Range("K134").Value = 5845.6155
$ws.Range("M134").Value = -3310.6155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 701.0625
$ws.Range("I7").Value = 708.55554
$ws.Range("J7").Value = 691.4286
$ws.Range("K7").Value = 708.55554
$ws.Range("L7").Value = 691.4286
$ws.Range("M7").Value = -595.55554
$ws.Range("N7").Value = -917.4286
$ws.Range("H31").Value = 2559.8076
$ws.Range("I31").Value = 2065.4375
$ws.Range("J31").Value = 3350.8
$ws.Range("K31").Value = 2065.4375
$ws.Range("L31").Value = 3350.8
$ws.Range("M31").Value = -1770.4375
$ws.Range("N31").Value = -3940.8
$ws.Range("H34").Value = 2559.8076
$ws.Range("I34").Value = 2065.4375
$ws.Range("J34").Value = 3350.8
$ws.Range("K34").Value = 2065.4375
$ws.Range("L34").Value = 3350.8
$ws.Range("M34").Value = -1863.4375
$ws.Range("N34").Value = -3754.8
$ws.Range("H58").Value = 145570.86
$ws.Range("I58").Value = 169282.67
$ws.Range("J58").Value = 3300
$ws.Range("K58").Value = 169282.67
$ws.Range("L58").Value = 3300
$ws.Range("M58").Value = -169079.67
$ws.Range("N58").Value = -3706
$ws.Range("H136").Value = 145570.86
$ws.Range("I136").Value = 169282.67
$ws.Range("J136").Value = 3300
$ws.Range("K136").Value = 507848.01
$ws.Range("L136").Value = 9900
$ws.Range("M136").Value = -505298.01
$ws.Range("N136").Value = -15000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 280.66666
$ws.Range("J12").Value = 325.7143
$ws.Range("L12").Value = 977.1428999999999
$ws.Range("N12").Value = -1323.1429
$ws.Range("H129").Value = 3012.8
$ws.Range("I129").Value = 2792
$ws.Range("K129").Value = 8376
$ws.Range("M129").Value = -3376
$ws.Range("H131").Value = 2179166.5
$ws.Range("I131").Value = 2553.6667
$ws.Range("K131").Value = 7661.000100000001
$ws.Range("M131").Value = -2621.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H136").Value = 42617.89
$ws.Range("J136").Value = 42617.89
$ws.Range("L136").Value = 127853.67
$ws.Range("N136").Value = -132953.67
$ws.Range("H141").Value = 50000
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 70000
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 70000
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -80360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6953.2144
$ws.Range("J46").Value = 4965.6113
$ws.Range("L46").Value = 4965.6113
$ws.Range("N46").Value = -5341.6113
$ws.Range("H82").Value = 2285.65
$ws.Range("I82").Value = 1138.2
$ws.Range("J82").Value = 2668.1333
$ws.Range("K82").Value = 1138.2
$ws.Range("L82").Value = 2668.1333
$ws.Range("M82").Value = -777.2
$ws.Range("N82").Value = -3390.1333
$ws.Range("H85").Value = 2285.65
$ws.Range("I85").Value = 1138.2
$ws.Range("J85").Value = 2668.1333
$ws.Range("K85").Value = 1138.2
$ws.Range("L85").Value = 2668.1333
$ws.Range("M85").Value = 109.8
$ws.Range("N85").Value = -5164.1333
$ws.Range("H136").Value = 3749.2856
$ws.Range("I136").Value = 2912
$ws.Range("K136").Value = 8736
$ws.Range("M136").Value = -6186

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 460.84616
$ws.Range("I107").Value = 476.9091
$ws.Range("J107").Value = 372.5
$ws.Range("K107").Value = 1430.7273
$ws.Range("L107").Value = 1117.5
$ws.Range("M107").Value = 489.2727
$ws.Range("N107").Value = -4957.5
$ws.Range("H136").Value = 2344.1482
$ws.Range("I136").Value = 1973.409
$ws.Range("J136").Value = 3975.4
$ws.Range("K136").Value = 5920.227000000001
$ws.Range("L136").Value = 11926.2
$ws.Range("M136").Value = -3370.227000000001
$ws.Range("N136").Value = -17026.2
$ws.Range("H137").Value = 79900
$ws.Range("J137").Value = 79900
$ws.Range("L137").Value = 79900
$ws.Range("N137").Value = -90100
$ws.Range("H141").Value = 74999.5
$ws.Range("J141").Value = 74999.5
$ws.Range("L141").Value = 74999.5
$ws.Range("N141").Value = -85359.5

Write-Host "Applied market price updates across 8 sheets"